$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2:C5 store numeric-looking results as text (matches the original
# inline-string cell type), so force Text format before assigning the
# new values to avoid Excel's automatic numeric conversion.
$ws.Range("B2:C5").NumberFormat = "@"

$ws.Range("B2").Value = "0.057947338"
$ws.Range("C2").Value = "0.12490211"

$ws.Range("B3").Value = "0.04426039"
$ws.Range("C3").Value = "0.13343696"

$ws.Range("B4").Value = "0.10513334"
$ws.Range("C4").Value = "0.26206714"

$ws.Range("B5").Value = "0.07529135"
$ws.Range("C5").Value = "0.21334727"
